$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1743614.6
$ws.Range("I6").Value = 7143182
$ws.Range("J6").Value = 303730
$ws.Range("K6").Value = 21429546
$ws.Range("L6").Value = 911190
$ws.Range("M6").Value = -21429434
$ws.Range("N6").Value = -911414
$ws.Range("H132").Value = 7872.5405
$ws.Range("I132").Value = 6871.6
$ws.Range("J132").Value = 9957.833000000001
$ws.Range("K132").Value = 20614.8
$ws.Range("L132").Value = 29873.499
$ws.Range("M132").Value = -18084.8
$ws.Range("N132").Value = -34933.499
$ws.Range("H138").Value = 2057.3447
$ws.Range("I138").Value = 2009.8948
$ws.Range("K138").Value = 6029.6844
$ws.Range("M138").Value = -889.6844000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7638.8486
$ws.Range("I32").Value = 6859.746
$ws.Range("K32").Value = 6859.746
$ws.Range("M32").Value = -6572.746
$ws.Range("H110").Value = 1670.7059
$ws.Range("I110").Value = 1463.4166
$ws.Range("J110").Value = 2168.2
$ws.Range("K110").Value = 1463.4166
$ws.Range("L110").Value = 2168.2
$ws.Range("M110").Value = 581.5834
$ws.Range("N110").Value = -6258.2
$ws.Range("H122").Value = 1504.2632
$ws.Range("I122").Value = 1256.3125
$ws.Range("K122").Value = 3768.9375
$ws.Range("M122").Value = -1318.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2203.6924
$ws.Range("I105").Value = 1997.75
$ws.Range("J105").Value = 2295.2222
$ws.Range("K105").Value = 1997.75
$ws.Range("L105").Value = 2295.2222
$ws.Range("M105").Value = -250.75
$ws.Range("N105").Value = -5789.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3137.3845
$ws.Range("I16").Value = 3208.6
$ws.Range("J16").Value = 2900
$ws.Range("K16").Value = 3208.6
$ws.Range("L16").Value = 2900
$ws.Range("M16").Value = -2921.6
$ws.Range("N16").Value = -3474
$ws.Range("H31").Value = 2973.78
$ws.Range("I31").Value = 2843.889
$ws.Range("J31").Value = 3046.8438
$ws.Range("K31").Value = 2843.889
$ws.Range("L31").Value = 3046.8438
$ws.Range("M31").Value = -2548.889
$ws.Range("N31").Value = -3636.8438
$ws.Range("H34").Value = 2973.78
$ws.Range("I34").Value = 2843.889
$ws.Range("J34").Value = 3046.8438
$ws.Range("K34").Value = 2843.889
$ws.Range("L34").Value = 3046.8438
$ws.Range("M34").Value = -2641.889
$ws.Range("N34").Value = -3450.8438
$ws.Range("H105").Value = 887
$ws.Range("I105").Value = 902.94446
$ws.Range("J105").Value = 600
$ws.Range("K105").Value = 902.94446
$ws.Range("L105").Value = 600
$ws.Range("M105").Value = 844.05554
$ws.Range("N105").Value = -4094
$ws.Range("H113").Value = 3137.3845
$ws.Range("I113").Value = 3208.6
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 3208.6
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -1038.6
$ws.Range("N113").Value = -7240
$ws.Range("H134").Value = 2559.6
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 2599.5715
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 7798.7145
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -12868.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 580.2
$ws.Range("I7").Value = 201
$ws.Range("J7").Value = 675
$ws.Range("K7").Value = 603
$ws.Range("L7").Value = 2025
$ws.Range("M7").Value = -491
$ws.Range("N7").Value = -2249
$ws.Range("H80").Value = 1598
$ws.Range("I80").Value = 1400
$ws.Range("J80").Value = 1677.2
$ws.Range("K80").Value = 4200
$ws.Range("L80").Value = 5031.6
$ws.Range("M80").Value = -3264
$ws.Range("N80").Value = -6903.6
$ws.Range("H83").Value = 1598
$ws.Range("I83").Value = 1400
$ws.Range("J83").Value = 1677.2
$ws.Range("K83").Value = 12600
$ws.Range("L83").Value = 15094.8
$ws.Range("M83").Value = -7920
$ws.Range("N83").Value = -24454.8
$ws.Range("H92").Value = 1211.3334
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1211.3334
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 3634.0002
$ws.Range("N92").Value = -6130.0002
$ws.Range("M92").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8291.651
$ws.Range("I102").Value = 7411.6772
$ws.Range("J102").Value = 10564.917
$ws.Range("K102").Value = 7411.6772
$ws.Range("L102").Value = 10564.917
$ws.Range("M102").Value = -5789.6772
$ws.Range("N102").Value = -13808.917
$ws.Range("H122").Value = 3563.5417
$ws.Range("I122").Value = 3226.1875
$ws.Range("J122").Value = 4238.25
$ws.Range("K122").Value = 9678.5625
$ws.Range("L122").Value = 12714.75
$ws.Range("M122").Value = -7228.5625
$ws.Range("N122").Value = -17614.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2596.1333
$ws.Range("I7").Value = 2234.2
$ws.Range("J7").Value = 3320
$ws.Range("K7").Value = 2234.2
$ws.Range("L7").Value = 3320
$ws.Range("M7").Value = -2122.2
$ws.Range("N7").Value = -3544
$ws.Range("H46").Value = 556073.9
$ws.Range("I46").Value = 522.5
$ws.Range("J46").Value = 714802.9
$ws.Range("K46").Value = 522.5
$ws.Range("L46").Value = 714802.9
$ws.Range("M46").Value = -334.5
$ws.Range("N46").Value = -715178.9
$ws.Range("H126").Value = 2596.1333
$ws.Range("I126").Value = 2234.2
$ws.Range("J126").Value = 3320
$ws.Range("K126").Value = 6702.599999999999
$ws.Range("L126").Value = 9960
$ws.Range("M126").Value = -4232.599999999999
$ws.Range("N126").Value = -14900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6073.3335
$ws.Range("I62").Value = 4666
$ws.Range("J62").Value = 8888
$ws.Range("K62").Value = 4666
$ws.Range("L62").Value = 8888
$ws.Range("M62").Value = -4042
$ws.Range("N62").Value = -10136
$ws.Range("H65").Value = 6073.3335
$ws.Range("I65").Value = 4666
$ws.Range("J65").Value = 8888
$ws.Range("K65").Value = 23330
$ws.Range("L65").Value = 44440
$ws.Range("M65").Value = -20210
$ws.Range("N65").Value = -50680
$ws.Range("H96").Value = 2994.6667
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2994.6667
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2994.6667
$ws.Range("N96").Value = -5740.6667
$ws.Range("H122").Value = 1964.5714
$ws.Range("I122").Value = 1273.9412
$ws.Range("J122").Value = 4899.75
$ws.Range("K122").Value = 3821.8236
$ws.Range("L122").Value = 14699.25
$ws.Range("M122").Value = -1371.8236
$ws.Range("N122").Value = -19599.25
$ws.Range("H132").Value = 1975.75
$ws.Range("I132").Value = 1696.7931
$ws.Range("J132").Value = 2275.3704
$ws.Range("K132").Value = 5090.379300000001
$ws.Range("L132").Value = 6826.111199999999
$ws.Range("M132").Value = -2560.379300000001
$ws.Range("N132").Value = -11886.1112
$ws.Range("H136").Value = 1530.8103
$ws.Range("I136").Value = 1247.4615
$ws.Range("K136").Value = 3742.3845
$ws.Range("M136").Value = -1192.3845
$ws.Range("M96").ClearContents()
